$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.264.96'
$ws.Range('E2').Value = '  -1.79%  '
$ws.Range('D3').Value = '2.274.98'
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'298.27"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.76%  '
$ws.Range('D6').Value = "'95.52"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.54%  '
$ws.Range('E8').Value = '  -3.80%  '
$ws.Range('D9').Value = "'0.493"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.64%  '
$ws.Range('D10').Value = "'33.37"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.49%  '
$ws.Range('D11').Value = "'0.0789"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.06%  '
$ws.Range('D12').Value = "'48.49"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -7.66%  '
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('E14').Value = '  -3.26%  '
$ws.Range('D15').Value = "'15.72"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('D16').Value = '2.628.81'
$ws.Range('E16').Value = '  -2.83%  '
$ws.Range('D17').Value = '2.282.92'
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('D18').Value = "'0.781"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.70%  '
$ws.Range('D19').Value = '42.219.41'
$ws.Range('E19').Value = '  -1.70%  '
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').Value = "'11.53"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0891'
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('D22').Value = "'6.01"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.74%  '
$ws.Range('D23').Value = "'66.69"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('D24').Value = "'233.52"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.36%  '
$ws.Range('D25').Value = "'1.98"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.72%  '
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('E27').Value = '  -4.37%  '
$ws.Range('D28').Value = "'23.93"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.98%  '
$ws.Range('E29').Value = '  -1.08%  '
$ws.Range('D30').Value = "'167.22"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.11%  '
$ws.Range('D31').Value = "'34.04"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.65%  '
$ws.Range('D32').Value = "'9.07"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.61%  '
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('E34').Value = '  -4.02%  '
$ws.Range('D35').Value = "'4.53"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.30%  '
$ws.Range('D36').Value = "'0.0692"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.02%  '
$ws.Range('E37').Value = '  -5.00%  '
$ws.Range('D38').Value = "'16.31"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -7.39%  '
$ws.Range('E39').Value = '  -4.21%  '
$ws.Range('E40').Value = '  -3.06%  '
$ws.Range('E41').Value = '  -3.32%  '
$ws.Range('E42').Value = '  -7.08%  '
$ws.Range('D43').Value = "'2.39"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.77%  '
$ws.Range('D44').Value = '1.962.73'
$ws.Range('E44').Value = '  -3.30%  '
$ws.Range('D45').Value = "'0.0279"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('D46').Value = "'17.51"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.58%  '
$ws.Range('D47').Value = "'9.62"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.75%  '
$ws.Range('D48').Value = "'2.80"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.66%  '
$ws.Range('D49').Value = '2.500.06'
$ws.Range('E49').Value = '  -2.26%  '
$ws.Range('D50').Value = "'52.36"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.88%  '
$ws.Range('E51').Value = '  -5.43%  '
